$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column F slightly (stored width 23 -> 24)
# Excel's ColumnWidth property uses character units that get padded/rounded
# when serialized to the OOXML "width" attribute; 23.17 round-trips to 24.
$ws.Columns.Item(6).ColumnWidth = 23.17

# Row 3 (PORCELANATO) updated sales figures
$ws.Range("D3").Value = 3237.89
$ws.Range("E3").Value = 12233.6693
$ws.Range("F3").Value = 0.2092801337742344

# Row 4 (TOTAL) recalculated figures
$ws.Range("D4").Value = 74224.11
$ws.Range("E4").Value = -58752.5507
$ws.Range("F4").Value = 4.797455030922449
